$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Weekly refresh: prepend two new observations (Black Amber, Primera /
# Segunda, Región de O'Higgins) above the existing Ciruela block, pushing
# rows 122-166 down to 124-168.
$ws.Rows.Item(122).Resize(2).Insert()

$newRows = @(
    @(122, @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44588, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 26, 260000, 265000, 262500, "`$/bins (450 kilos)", "Región de O'Higgins", 583, 450)),
    @(123, @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44588, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Segunda", 20, 210000, 215000, 212500, "`$/bins (450 kilos)", "Región de O'Higgins", 472, 450))
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# Match the date-number-format style used throughout column D for the two
# newly written date cells (Insert already seeded s="2" on D122/D123, but
# set it explicitly too so the value sticks with the right format).
$ws.Range("D122:D123").NumberFormat = $ws.Range("D124").NumberFormat
